$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 529.625
$ws.Range("I28").Value = 176.71428
$ws.Range("K28").Value = 176.71428
$ws.Range("M28").Value = 308.28572

$ws.Range("H70").Value = 1695.5714
$ws.Range("I70").Value = 1820
$ws.Range("J70").Value = 1529.6666
$ws.Range("K70").Value = 5460
$ws.Range("L70").Value = 4588.9998
$ws.Range("M70").Value = -5190
$ws.Range("N70").Value = -5128.9998

$ws.Range("H73").Value = 1695.5714
$ws.Range("I73").Value = 1820
$ws.Range("J73").Value = 1529.6666
$ws.Range("K73").Value = 5460
$ws.Range("L73").Value = 4588.9998
$ws.Range("M73").Value = -4524
$ws.Range("N73").Value = -6460.9998

$ws.Range("H138").Value = 3298.2
$ws.Range("I138").Value = 2532.0908
$ws.Range("J138").Value = 3489.7273
$ws.Range("K138").Value = 7596.2724
$ws.Range("L138").Value = 10469.1819
$ws.Range("M138").Value = -2456.2724
$ws.Range("N138").Value = -20749.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 993.25714
$ws.Range("I2").Value = 1077.6207
$ws.Range("J2").Value = 585.5
$ws.Range("K2").Value = 1077.6207
$ws.Range("L2").Value = 585.5
$ws.Range("M2").Value = -964.6206999999999
$ws.Range("N2").Value = -811.5

$ws.Range("H45").Value = 1957.1428
$ws.Range("I45").Value = 1900
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 1900
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -1523
$ws.Range("N45").Value = -2754

$ws.Range("H116").Value = 993.25714
$ws.Range("I116").Value = 1077.6207
$ws.Range("J116").Value = 585.5
$ws.Range("K116").Value = 1077.6207
$ws.Range("L116").Value = 585.5
$ws.Range("M116").Value = 1216.3793
$ws.Range("N116").Value = -5173.5

$ws.Range("H132").Value = 1321235.9
$ws.Range("I132").Value = 1568161.5
$ws.Range("K132").Value = 4704484.5
$ws.Range("M132").Value = -4701954.5

$ws.Range("H133").Value = 69750
$ws.Range("J133").Value = 69750
$ws.Range("L133").Value = 69750
$ws.Range("N133").Value = -74810

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 993.25714
$ws.Range("I3").Value = 1077.6207
$ws.Range("J3").Value = 585.5
$ws.Range("K3").Value = 1077.6207
$ws.Range("L3").Value = 585.5
$ws.Range("M3").Value = -963.6206999999999
$ws.Range("N3").Value = -813.5

$ws.Range("H11").Value = 215.77777
$ws.Range("I11").Value = 82
$ws.Range("J11").Value = 282.66666
$ws.Range("K11").Value = 82
$ws.Range("L11").Value = 282.66666
$ws.Range("M11").Value = 58
$ws.Range("N11").Value = -562.66666

$ws.Range("H132").Value = 54997.5
$ws.Range("J132").Value = 54997.5
$ws.Range("L132").Value = 54997.5
$ws.Range("N132").Value = -65117.5

$ws.Range("I134").Value = 11755.917
$ws.Range("K134").Value = 35267.751
$ws.Range("M134").Value = -32732.751

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 12682.5
$ws.Range("I103").Value = 12682.5
$ws.Range("K103").Value = 12682.5
$ws.Range("M103").Value = -11510.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2631.1333
$ws.Range("J113").Value = 2837.2307
$ws.Range("L113").Value = 8511.6921
$ws.Range("N113").Value = -12851.6921

$ws.Range("H122").Value = 1614192.8
$ws.Range("J122").Value = 1579
$ws.Range("L122").Value = 14211
$ws.Range("N122").Value = -19111

$ws.Range("H126").Value = 14583.333
$ws.Range("J126").Value = 14583.333
$ws.Range("L126").Value = 43749.999
$ws.Range("N126").Value = -53629.999

$ws.Range("H129").Value = 5159.8
$ws.Range("J129").Value = 4866
$ws.Range("L129").Value = 14598
$ws.Range("N129").Value = -24598

$ws.Range("H130").Value = 12833.333
$ws.Range("I130").Value = 2000
$ws.Range("K130").Value = 6000
$ws.Range("M130").Value = -980

$ws.Range("H131").Value = 5559.087
$ws.Range("I131").Value = 1731.1111
$ws.Range("K131").Value = 5193.3333
$ws.Range("M131").Value = -153.3333000000002

$ws.Range("H137").Value = 10554.059
$ws.Range("I137").Value = 3853.3333
$ws.Range("J137").Value = 14209
$ws.Range("K137").Value = 11559.9999
$ws.Range("L137").Value = 42627
$ws.Range("M137").Value = -6459.999899999999
$ws.Range("N137").Value = -52827

$ws.Range("H138").Value = 80257.5
$ws.Range("J138").Value = 15000
$ws.Range("L138").Value = 45000
$ws.Range("N138").Value = -55280

$ws.Range("H140").Value = 2170.08
$ws.Range("I140").Value = 1923.9131
$ws.Range("J140").Value = 5001
$ws.Range("K140").Value = 5771.7393
$ws.Range("L140").Value = 15003
$ws.Range("M140").Value = -591.7393000000002
$ws.Range("N140").Value = -25363

$ws.Range("H141").Value = 2615.3333
$ws.Range("I141").Value = 2615.3333
$ws.Range("K141").Value = 7845.999899999999
$ws.Range("M141").Value = -2665.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 77094.07000000001
$ws.Range("I122").Value = 82609
$ws.Range("K122").Value = 247827
$ws.Range("M122").Value = -245377

$ws.Range("H132").Value = 24298.5
$ws.Range("I132").Value = 22373.188
$ws.Range("K132").Value = 67119.564
$ws.Range("M132").Value = -64589.564

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 53436.367
$ws.Range("I16").Value = 63255.125
$ws.Range("J16").Value = 1069.6666
$ws.Range("K16").Value = 63255.125
$ws.Range("L16").Value = 1069.6666
$ws.Range("M16").Value = -63085.125
$ws.Range("N16").Value = -1409.6666

$ws.Range("H132").Value = 4497724
$ws.Range("I132").Value = 8991646
$ws.Range("J132").Value = 3801.923
$ws.Range("K132").Value = 26974938
$ws.Range("L132").Value = 11405.769
$ws.Range("M132").Value = -26972408
$ws.Range("N132").Value = -16465.769

$ws.Range("H133").Value = 74979.44500000001
$ws.Range("J133").Value = 74979.44500000001
$ws.Range("L133").Value = 74979.44500000001
$ws.Range("N133").Value = -80039.44500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 36090.6
$ws.Range("I107").Value = 6925
$ws.Range("J107").Value = 48590.145
$ws.Range("K107").Value = 20775
$ws.Range("L107").Value = 145770.435
$ws.Range("M107").Value = -18855
$ws.Range("N107").Value = -149610.435

$ws.Range("H132").Value = 3970013.8
$ws.Range("I132").Value = 4387639
$ws.Range("K132").Value = 13162917
$ws.Range("M132").Value = -13160387
